$wb = $excel.ActiveWorkbook

# Sheet ALC - Row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 80.388885
$ws.Range("I33").Value = 73.5
$ws.Range("K33").Value = 73.5
$ws.Range("M33").Value = 155.5

# Sheet ALC - Row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1600
$ws.Range("J43").Value = 1600
$ws.Range("L43").Value = 1600
$ws.Range("N43").Value = -1738

# Sheet ALC - Row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3943.6667
$ws.Range("J106").Value = 5469.3335
$ws.Range("L106").Value = 5469.3335
$ws.Range("N106").Value = -6731.3335

# Sheet ALC - Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1007.5833
$ws.Range("I132").Value = 932.4386
$ws.Range("K132").Value = 2797.3158
$ws.Range("M132").Value = -267.3157999999999

# Sheet ALC - Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4997.379
$ws.Range("I138").Value = 6237.385
$ws.Range("J138").Value = 4693.2266
$ws.Range("K138").Value = 18712.155
$ws.Range("L138").Value = 14079.6798
$ws.Range("M138").Value = -13572.155
$ws.Range("N138").Value = -24359.6798

# Sheet ARM - Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11161.718
$ws.Range("I32").Value = 9730.777
$ws.Range("K32").Value = 9730.777
$ws.Range("M32").Value = -9443.777

# Sheet ARM - Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 19987.209
$ws.Range("I61").Value = 25954.576
$ws.Range("K61").Value = 25954.576
$ws.Range("M61").Value = -25742.576

# Sheet ARM - Row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2199.818
$ws.Range("I97").Value = 1650.5
$ws.Range("K97").Value = 1650.5
$ws.Range("M97").Value = -1154.5

# Sheet ARM - Row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 954.56665
$ws.Range("I102").Value = 946.7692
$ws.Range("J102").Value = 1005.25
$ws.Range("K102").Value = 946.7692
$ws.Range("L102").Value = 1005.25
$ws.Range("M102").Value = 675.2308
$ws.Range("N102").Value = -4249.25

# Sheet ARM - Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 19987.209
$ws.Range("I136").Value = 25954.576
$ws.Range("K136").Value = 77863.728
$ws.Range("M136").Value = -75313.728

# Sheet ARM - Row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 41143.855
$ws.Range("J139").Value = 41143.855
$ws.Range("L139").Value = 41143.855
$ws.Range("N139").Value = -51423.855

# Sheet BSM - Row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2562.7368
$ws.Range("I20").Value = 2275.5833
$ws.Range("K20").Value = 2275.5833
$ws.Range("M20").Value = -2028.5833

# Sheet BSM - Row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2323.745
$ws.Range("I105").Value = 2289.5952
$ws.Range("K105").Value = 2289.5952
$ws.Range("M105").Value = -542.5952000000002

# Sheet CRP - Row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1566.4445
$ws.Range("I22").Value = 766.6667
$ws.Range("J22").Value = 1966.3334
$ws.Range("K22").Value = 766.6667
$ws.Range("L22").Value = 1966.3334
$ws.Range("M22").Value = -416.6667
$ws.Range("N22").Value = -2666.3334

# Sheet CRP - Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2477.182
$ws.Range("I31").Value = 2125
$ws.Range("K31").Value = 2125
$ws.Range("M31").Value = -1830

# Sheet CRP - Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2477.182
$ws.Range("I34").Value = 2125
$ws.Range("K34").Value = 2125
$ws.Range("M34").Value = -1923

# Sheet CRP - Row 70
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 50045
$ws.Range("J70").Value = 50045
$ws.Range("L70").Value = 50045
$ws.Range("N70").Value = -50675

# Sheet CRP - Row 73
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 50045
$ws.Range("J73").Value = 50045
$ws.Range("L73").Value = 50045
$ws.Range("N73").Value = -52229

# Sheet CRP - Row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# Sheet CUL - Row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 429
$ws.Range("I2").Value = 286.25
$ws.Range("K2").Value = 1717.5
$ws.Range("M2").Value = -1604.5

# Sheet CUL - Row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 466.42856
$ws.Range("I18").Value = 247
$ws.Range("K18").Value = 741
$ws.Range("M18").Value = -572

# Sheet CUL - Row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 4736.6
$ws.Range("J103").Value = 5218
$ws.Range("L103").Value = 15654
$ws.Range("N103").Value = -17412

# Sheet CUL - Row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1939.4
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1939.4
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 17454.6
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -22354.6

# Sheet CUL - Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11842.306
$ws.Range("J131").Value = 12333.131
$ws.Range("L131").Value = 36999.393
$ws.Range("N131").Value = -47079.393

# Sheet CUL - Row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2033.3334
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2033.3334
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 18300.0006
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -23360.0006

# Sheet GSM - Row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Sheet GSM - Row 62
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 39999.5
$ws.Range("J62").Value = 39999.5
$ws.Range("L62").Value = 39999.5
$ws.Range("N62").Value = -41371.5

# Sheet GSM - Row 65
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 39999.5
$ws.Range("J65").Value = 39999.5
$ws.Range("L65").Value = 119998.5
$ws.Range("N65").Value = -126862.5

# Sheet GSM - Row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1646.6666
$ws.Range("I80").Value = 980
$ws.Range("J80").Value = 1980
$ws.Range("K80").Value = 980
$ws.Range("L80").Value = 1980
$ws.Range("M80").Value = 18
$ws.Range("N80").Value = -3976

# Sheet GSM - Row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Sheet GSM - Row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 1646.6666
$ws.Range("I83").Value = 980
$ws.Range("J83").Value = 1980
$ws.Range("K83").Value = 4900
$ws.Range("L83").Value = 9900
$ws.Range("M83").Value = 92
$ws.Range("N83").Value = -19884

# Sheet GSM - Row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Sheet LTW - Row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3699.75
$ws.Range("I46").Value = 2900
$ws.Range("J46").Value = 3966.3333
$ws.Range("K46").Value = 2900
$ws.Range("L46").Value = 3966.3333
$ws.Range("M46").Value = -2900
$ws.Range("N46").Value = -4342.3333

